$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("J3").Value = 4.2
$ws.Range("X3").Value = 17
$ws.Range("Y3").Value = 29
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 14.5
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 13.5
$ws.Range("AG3").Value = 12.5
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 21
$ws.Range("AK3").Value = 21
$ws.Range("AL3").Value = 36
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 10.5
$ws.Range("AO3").Value = 1000
$ws.Range("G4").Value = 5.8
$ws.Range("H4").Value = 1.61
$ws.Range("Q4").Value = 1.69
$ws.Range("U4").Value = 2.22
$ws.Range("W4").Value = 1.21
$ws.Range("X4").Value = 21
$ws.Range("Y4").Value = 11
$ws.Range("AA4").Value = 15.5
$ws.Range("AI4").Value = 28
$ws.Range("AJ4").Value = 180
$ws.Range("F5").Value = 2.16
$ws.Range("G5").Value = 2.36
$ws.Range("I5").Value = 3.65
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 3.85
$ws.Range("V5").Value = 1.38
$ws.Range("W5").Value = 1.75
$ws.Range("N6").Value = 5.4
$ws.Range("S6").Value = 2.08
$ws.Range("AD6").Value = 21
$ws.Range("AF6").Value = 14.5
$ws.Range("AG6").Value = 11.5
$ws.Range("AN6").Value = 6.8
$ws.Range("F7").Value = 8.4
$ws.Range("G7").Value = 8.800000000000001
$ws.Range("J7").Value = 6.2
$ws.Range("K7").Value = 6.4
$ws.Range("T7").Value = 1.64
$ws.Range("U7").Value = 2.46
$ws.Range("AB7").Value = 48
$ws.Range("AF7").Value = 85
$ws.Range("AG7").Value = 32
$ws.Range("AJ7").Value = 310
$ws.Range("I8").Value = 2.72
$ws.Range("N9").Value = 1.01
$ws.Range("Q9").Value = 1.31
$ws.Range("S9").Value = 1.31
$ws.Range("F10").Value = 5.3
$ws.Range("J10").Value = 3.7
$ws.Range("K10").Value = 85
$ws.Range("Q10").Value = 1.76
$ws.Range("V10").Value = 2.46
$ws.Range("F11").Value = 1.49
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 4.8
$ws.Range("I11").Value = 7
$ws.Range("J11").Value = 3.35
$ws.Range("K11").Value = 6
$ws.Range("S11").Value = 2.16
$ws.Range("T11").Value = 1.01
$ws.Range("U11").Value = 1.01
$ws.Range("V11").Value = 1.16
$ws.Range("W11").Value = 2.34
$ws.Range("X11").Value = 30
$ws.Range("Z11").Value = 60
$ws.Range("AA11").Value = 1000
$ws.Range("AB11").Value = 970
$ws.Range("AC11").Value = 970
$ws.Range("AD11").Value = 27
$ws.Range("AE11").Value = 1000
$ws.Range("AF11").Value = 970
$ws.Range("AG11").Value = 970
$ws.Range("AH11").Value = 22
$ws.Range("AI11").Value = 1000
$ws.Range("AJ11").Value = 970
$ws.Range("AK11").Value = 970
$ws.Range("AL11").Value = 32
$ws.Range("AM11").Value = 1000
$ws.Range("AN11").Value = 970
$ws.Range("AO11").Value = 70
$ws.Range("P12").Value = 2.16
$ws.Range("AB12").Value = 970
$ws.Range("AC12").Value = 970
$ws.Range("AF12").Value = 970
$ws.Range("AG12").Value = 970
$ws.Range("AH12").Value = 970
$ws.Range("AK12").Value = 970
$ws.Range("AN12").Value = 970
$ws.Range("I13").Value = 2.34
$ws.Range("P13").Value = 2.08
$ws.Range("Q13").Value = 1.75
$ws.Range("V13").Value = 1.74
$ws.Range("X13").Value = 21
$ws.Range("AA13").Value = 30
$ws.Range("AC13").Value = 9.4
$ws.Range("AF13").Value = 30
$ws.Range("AM13").Value = 90
$ws.Range("F14").Value = 1.33
$ws.Range("G14").Value = 1.56
$ws.Range("H14").Value = 7.2
$ws.Range("I14").Value = 65
$ws.Range("J14").Value = 4.3
$ws.Range("K14").Value = 7.4
$ws.Range("R14").Value = 1.39
$ws.Range("S14").Value = 2.4
$ws.Range("T14").Value = 1.7
$ws.Range("U14").Value = 1.62
$ws.Range("V14").Value = 1.02
$ws.Range("W14").Value = 2.78
$ws.Range("F15").Value = 2.32
$ws.Range("H15").Value = 2.76
$ws.Range("I15").Value = 4.5
$ws.Range("J15").Value = 2.9
$ws.Range("P15").Value = 1.68
$ws.Range("H16").Value = 2.6
$ws.Range("I16").Value = 2.64
$ws.Range("Q16").Value = 1.83
$ws.Range("V16").Value = 1.61
$ws.Range("AA16").Value = 36
$ws.Range("AE16").Value = 26
$ws.Range("AO16").Value = 19
$ws.Range("G17").Value = 4.7
$ws.Range("T17").Value = 1.94
$ws.Range("W17").Value = 1.27
$ws.Range("X17").Value = 12.5
$ws.Range("AA17").Value = 22
$ws.Range("I18").Value = 2.68
$ws.Range("L18").Value = 1.44
$ws.Range("N18").Value = 3.85
$ws.Range("O18").Value = 1.34
$ws.Range("P18").Value = 1.93
$ws.Range("Q18").Value = 2.04
$ws.Range("R18").Value = 1.36
$ws.Range("V18").Value = 1.59
$ws.Range("AA18").Value = 38
$ws.Range("AB18").Value = 11.5
$ws.Range("AI18").Value = 42
$ws.Range("AM18").Value = 95
$ws.Range("AO18").Value = 25
$ws.Range("K19").Value = 4.4
$ws.Range("N19").Value = 5.5
$ws.Range("P19").Value = 2.5
$ws.Range("Q19").Value = 1.65
$ws.Range("X19").Value = 21
$ws.Range("AA19").Value = 20
$ws.Range("AF19").Value = 36
$ws.Range("AJ19").Value = 100
$ws.Range("AN19").Value = 38
$ws.Range("F20").Value = 1.6
$ws.Range("G20").Value = 1.62
$ws.Range("T20").Value = 1.58
$ws.Range("U20").Value = 2.6
$ws.Range("W20").Value = 2.6
$ws.Range("AN20").Value = 5.5
$ws.Range("Q21").Value = 1.39
$ws.Range("S21").Value = 1.97
$ws.Range("T21").Value = 1.73
$ws.Range("Y21").Value = 55
$ws.Range("AI21").Value = 95
$ws.Range("AL21").Value = 25
$ws.Range("AO21").Value = 110
$ws.Range("Q22").Value = 1.79
$ws.Range("S22").Value = 3
$ws.Range("X22").Value = 17.5
$ws.Range("AH22").Value = 15.5
$ws.Range("AN22").Value = 18.5
